# The author renamed the single worksheet from "Sheet1" to "Tables"
# (it now hosts the database "Tables" created for the app, per the commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "Tables"
